$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $cell = $ws.Range($range)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue "D2" "244.57"
Set-TextValue "D3" "21.92"
Set-TextValue "D4" "5.394"
Set-TextValue "D5" "0.05991"
Set-TextValue "D7" "0.8146"
Set-TextValue "D8" "0.9528"
Set-TextValue "B9" "One"
Set-TextValue "C9" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D9" "0.0005902"
Set-TextValue "E9" "8OneONE"
Set-TextValue "B10" "WazirX"
Set-TextValue "C10" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D10" "0.1425"
Set-TextValue "E10" "9WazirXWRX"
Set-TextValue "B11" "MandalaExchangeToken"
Set-TextValue "C11" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D11" "0.07437"
Set-TextValue "E11" "10MandalaExchangeTokenMDX"
Set-TextValue "B12" "LiechtensteinCryptoassetsExchange"
Set-TextValue "C12" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D12" "0.03275"
Set-TextValue "E12" "11LiechtensteinCryptoassetsExchangeLCX"
Set-TextValue "B13" "BitrueCoin"
Set-TextValue "C13" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D13" "0.03054"
Set-TextValue "E13" "12BitrueCoinBTR"
Set-TextValue "B14" "BitMartToken"
Set-TextValue "C14" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D14" "0.09412"
Set-TextValue "E14" "13BitMartTokenBMX"
Set-TextValue "B15" "MCDex"
Set-TextValue "C15" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D15" "4.003"
Set-TextValue "E15" "14MCDexMCB"
Set-TextValue "B16" "BitForexToken"
Set-TextValue "C16" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D16" "0.001596"
Set-TextValue "E16" "15BitForexTokenBF"
Set-TextValue "B17" "CoinExToken"
Set-TextValue "C17" "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D17" "0.04813"
Set-TextValue "E17" "16CoinExTokenCET"
Set-TextValue "D18" "0.005497"
Set-TextValue "D19" "0.004152"
Set-TextValue "D20" "0.0009895"
Set-TextValue "D22" "3.678"
Set-TextValue "D23" "6.436"
Set-TextValue "D24" "2.188"
Set-TextValue "B41" "KickToken"
Set-TextValue "C41" "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D41" "0.006472"
Set-TextValue "E41" "40KickTokenKICKBestin24h"
Set-TextValue "B42" "BKEXToken"
Set-TextValue "C42" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1073"
Set-TextValue "E42" "41BKEXTokenBKK"
Set-TextValue "B43" "CEJI"
Set-TextValue "C43" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.002901"
Set-TextValue "E43" "42CEJICEJI"
Set-TextValue "D44" "0.006302"
Set-TextValue "D45" "0.00005146"
Set-TextValue "D47" "0.8603"
Set-TextValue "E47" "46CoinbaseStockTokenCOIN"
Set-TextValue "D48" "0.006142"
